$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 46064
$ws.Range("C3").Value = 46064
$ws.Range("C4").Value = 46064
$ws.Range("A5").Value = "A 61558-2023"
$ws.Range("B5").Value = 45265
$ws.Range("C5").Value = 46064
$ws.Range("F5").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G5").Value = 1.5
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("R5").Value = "Ask"
$ws.Range("S5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/artfynd/A 61558-2023 artfynd.xlsx`", `"A 61558-2023`")"
$ws.Range("T5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/kartor/A 61558-2023 karta.png`", `"A 61558-2023`")"
$ws.Range("V5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomål/A 61558-2023 FSC-klagomål.docx`", `"A 61558-2023`")"
$ws.Range("W5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomålsmail/A 61558-2023 FSC-klagomål mail.docx`", `"A 61558-2023`")"
$ws.Range("X5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsyn/A 61558-2023 tillsynsbegäran.docx`", `"A 61558-2023`")"
$ws.Range("Y5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsynsmail/A 61558-2023 tillsynsbegäran mail.docx`", `"A 61558-2023`")"
$ws.Range("A6").Value = "A 13766-2023"
$ws.Range("B6").Value = 45007
$ws.Range("C6").Value = 46064
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = 0.9
$ws.Range("H6").Value = 1
$ws.Range("L6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("R6").Value = "Större vattensalamander"
$ws.Range("S6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/artfynd/A 13766-2023 artfynd.xlsx`", `"A 13766-2023`")"
$ws.Range("T6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/kartor/A 13766-2023 karta.png`", `"A 13766-2023`")"
$ws.Range("V6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomål/A 13766-2023 FSC-klagomål.docx`", `"A 13766-2023`")"
$ws.Range("W6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomålsmail/A 13766-2023 FSC-klagomål mail.docx`", `"A 13766-2023`")"
$ws.Range("X6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsyn/A 13766-2023 tillsynsbegäran.docx`", `"A 13766-2023`")"
$ws.Range("Y6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsynsmail/A 13766-2023 tillsynsbegäran mail.docx`", `"A 13766-2023`")"
$ws.Range("A7").Value = "A 60891-2024"
$ws.Range("B7").Value = 45644
$ws.Range("C7").Value = 46064
$ws.Range("G7").Value = 16.1
$ws.Range("H7").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("R7").Value = "Lövgroda"
$ws.Range("S7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/artfynd/A 60891-2024 artfynd.xlsx`", `"A 60891-2024`")"
$ws.Range("T7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/kartor/A 60891-2024 karta.png`", `"A 60891-2024`")"
$ws.Range("V7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomål/A 60891-2024 FSC-klagomål.docx`", `"A 60891-2024`")"
$ws.Range("W7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomålsmail/A 60891-2024 FSC-klagomål mail.docx`", `"A 60891-2024`")"
$ws.Range("X7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsyn/A 60891-2024 tillsynsbegäran.docx`", `"A 60891-2024`")"
$ws.Range("Y7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsynsmail/A 60891-2024 tillsynsbegäran mail.docx`", `"A 60891-2024`")"
$ws.Range("C8").Value = 46064
$ws.Range("A9").Value = "A 31213-2023"
$ws.Range("B9").Value = 45113
$ws.Range("C9").Value = 46064
$ws.Range("G9").Value = 6.5
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 1
$ws.Range("O9").Value = 1
$ws.Range("R9").Value = "Skogsveronika"
$ws.Range("S9").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/artfynd/A 31213-2023 artfynd.xlsx`", `"A 31213-2023`")"
$ws.Range("T9").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/kartor/A 31213-2023 karta.png`", `"A 31213-2023`")"
$ws.Range("V9").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomål/A 31213-2023 FSC-klagomål.docx`", `"A 31213-2023`")"
$ws.Range("W9").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomålsmail/A 31213-2023 FSC-klagomål mail.docx`", `"A 31213-2023`")"
$ws.Range("X9").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsyn/A 31213-2023 tillsynsbegäran.docx`", `"A 31213-2023`")"
$ws.Range("Y9").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsynsmail/A 31213-2023 tillsynsbegäran mail.docx`", `"A 31213-2023`")"
$ws.Range("A10").Value = "A 34341-2024"
$ws.Range("B10").Value = 45525
$ws.Range("C10").Value = 46064
$ws.Range("F10").Value = "Övriga Aktiebolag"
$ws.Range("G10").Value = 14.4
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 1
$ws.Range("O10").Value = 1
$ws.Range("R10").Value = "Desmeknopp"
$ws.Range("S10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/artfynd/A 34341-2024 artfynd.xlsx`", `"A 34341-2024`")"
$ws.Range("T10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/kartor/A 34341-2024 karta.png`", `"A 34341-2024`")"
$ws.Range("V10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomål/A 34341-2024 FSC-klagomål.docx`", `"A 34341-2024`")"
$ws.Range("W10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/klagomålsmail/A 34341-2024 FSC-klagomål mail.docx`", `"A 34341-2024`")"
$ws.Range("X10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsyn/A 34341-2024 tillsynsbegäran.docx`", `"A 34341-2024`")"
$ws.Range("Y10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1286/tillsynsmail/A 34341-2024 tillsynsbegäran mail.docx`", `"A 34341-2024`")"
$ws.Range("C11").Value = 46064
$ws.Range("C12").Value = 46064
$ws.Range("C13").Value = 46064
$ws.Range("A14").Value = "A 34302-2024"
$ws.Range("B14").Value = 45524
$ws.Range("C14").Value = 46064
$ws.Range("G14").Value = 1.9
$ws.Range("A15").Value = "A 49536-2025"
$ws.Range("B15").Value = 45939.4221875
$ws.Range("C15").Value = 46064
$ws.Range("G15").Value = 1.5
$ws.Range("A16").Value = "A 32596-2024"
$ws.Range("B16").Value = 45513.61667824074
$ws.Range("C16").Value = 46064
$ws.Range("G16").Value = 2.6
$ws.Range("A17").Value = "A 49543-2025"
$ws.Range("B17").Value = 45939.42862268518
$ws.Range("C17").Value = 46064
$ws.Range("G17").Value = 1.4
$ws.Range("A18").Value = "A 49549-2025"
$ws.Range("B18").Value = 45939
$ws.Range("C18").Value = 46064
$ws.Range("G18").Value = 0.5
$ws.Range("A19").Value = "A 40417-2022"
$ws.Range("B19").Value = 44823
$ws.Range("C19").Value = 46064
$ws.Range("G19").Value = 2.3
$ws.Range("A20").Value = "A 18090-2022"
$ws.Range("B20").Value = 44684
$ws.Range("C20").Value = 46064
$ws.Range("G20").Value = 4.9
$ws.Range("A21").Value = "A 60803-2023"
$ws.Range("B21").Value = 45260
$ws.Range("C21").Value = 46064
$ws.Range("G21").Value = 1.6
$ws.Range("A22").Value = "A 38631-2023"
$ws.Range("B22").Value = 45162
$ws.Range("C22").Value = 46064
$ws.Range("G22").Value = 0.8
$ws.Range("C23").Value = 46064
